# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "302.82") need to stay
# TEXT (matching the source inline-string cells), otherwise Excel auto-converts
# them to numeric values and mangles formatting (trailing zeros, float noise).
# Pre-set an explicit text format on just those cells before writing the values.
$textCells = @(@(4,4), @(5,4), @(6,4), @(7,4), @(9,4), @(10,4), @(13,4), @(14,4), @(17,4), @(19,4), @(22,4), @(23,4), @(26,4), @(27,4), @(29,4), @(30,4), @(31,4), @(32,4), @(33,4), @(34,4), @(37,4), @(39,4), @(40,4), @(42,4), @(44,4), @(45,4), @(46,4), @(47,4), @(48,4))
foreach ($rc in $textCells) {
    $ws.Cells.Item($rc[0], $rc[1]).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '43.158.80'
$ws.Cells.Item(2, 5).Value = '  +2.42%  '
$ws.Cells.Item(3, 4).Value = '2.312.10'
$ws.Cells.Item(3, 5).Value = '  +1.76%  '
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '
$ws.Cells.Item(5, 4).Value = '302.82'
$ws.Cells.Item(5, 5).Value = '  +1.54%  '
$ws.Cells.Item(6, 4).Value = '100.27'
$ws.Cells.Item(6, 5).Value = '  +6.00%  '
$ws.Cells.Item(7, 4).Value = '0.506'
$ws.Cells.Item(7, 5).Value = '  +2.81%  '
$ws.Cells.Item(8, 5).Value = '  -0.11%  '
$ws.Cells.Item(9, 4).Value = '0.509'
$ws.Cells.Item(9, 5).Value = '  +3.91%  '
$ws.Cells.Item(10, 4).Value = '34.55'
$ws.Cells.Item(10, 5).Value = '  +4.34%  '
$ws.Cells.Item(11, 5).Value = '  +1.13%  '
$ws.Cells.Item(12, 5).Value = '  +4.12%  '
$ws.Cells.Item(13, 4).Value = '18.01'
$ws.Cells.Item(13, 5).Value = '  +14.76%  '
$ws.Cells.Item(14, 4).Value = '6.85'
$ws.Cells.Item(14, 5).Value = '  +3.40%  '
$ws.Cells.Item(15, 4).Value = '2.670.35'
$ws.Cells.Item(15, 5).Value = '  +1.70%  '
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '2.304.74'
$ws.Cells.Item(16, 5).Value = '  +1.59%  '
$ws.Cells.Item(17, 2).Value = 'Polygon'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(17, 4).Value = '0.820'
$ws.Cells.Item(17, 5).Value = '  +5.84%  '
$ws.Cells.Item(18, 4).Value = '43.102.77'
$ws.Cells.Item(18, 5).Value = '  +2.27%  '
$ws.Cells.Item(19, 4).Value = '12.62'
$ws.Cells.Item(19, 5).Value = '  +11.05%  '
$ws.Cells.Item(20, 5).Value = '  +2.20%  '
$ws.Cells.Item(21, 5).Value = '  +2.38%  '
$ws.Cells.Item(22, 4).Value = '67.88'
$ws.Cells.Item(22, 5).Value = '  +1.97%  '
$ws.Cells.Item(23, 4).Value = '237.95'
$ws.Cells.Item(23, 5).Value = '  +2.35%  '
$ws.Cells.Item(24, 5).Value = '  +13.35%  '
$ws.Cells.Item(25, 5).Value = '  +0.94%  '
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  +0.00%  '
$ws.Cells.Item(27, 4).Value = '24.86'
$ws.Cells.Item(27, 5).Value = '  +4.55%  '
$ws.Cells.Item(28, 5).Value = '  -4.88%  '
$ws.Cells.Item(29, 4).Value = '168.04'
$ws.Cells.Item(29, 5).Value = '  +1.03%  '
$ws.Cells.Item(30, 4).Value = '34.24'
$ws.Cells.Item(30, 5).Value = '  +1.88%  '
$ws.Cells.Item(31, 4).Value = '9.19'
$ws.Cells.Item(31, 5).Value = '  +1.85%  '
$ws.Cells.Item(32, 4).Value = '0.998'
$ws.Cells.Item(32, 5).Value = '  -0.12%  '
$ws.Cells.Item(33, 4).Value = '5.05'
$ws.Cells.Item(33, 5).Value = '  +3.09%  '
$ws.Cells.Item(34, 4).Value = '4.73'
$ws.Cells.Item(34, 5).Value = '  +5.74%  '
$ws.Cells.Item(35, 5).Value = '  +5.10%  '
$ws.Cells.Item(36, 5).Value = '  +6.60%  '
$ws.Cells.Item(37, 4).Value = '0.0693'
$ws.Cells.Item(37, 5).Value = '  +0.64%  '
$ws.Cells.Item(38, 5).Value = '  +4.16%  '
$ws.Cells.Item(39, 4).Value = '1.81'
$ws.Cells.Item(39, 5).Value = '  +5.08%  '
$ws.Cells.Item(40, 4).Value = '2.83'
$ws.Cells.Item(40, 5).Value = '  +2.04%  '
$ws.Cells.Item(41, 5).Value = '  +1.01%  '
$ws.Cells.Item(42, 4).Value = '2.35'
$ws.Cells.Item(42, 5).Value = '  -2.37%  '
$ws.Cells.Item(43, 4).Value = '2.001.39'
$ws.Cells.Item(43, 5).Value = '  +1.89%  '
$ws.Cells.Item(44, 4).Value = '0.0288'
$ws.Cells.Item(44, 5).Value = '  +4.02%  '
$ws.Cells.Item(45, 4).Value = '10.11'
$ws.Cells.Item(45, 5).Value = '  +5.78%  '
$ws.Cells.Item(46, 4).Value = '17.71'
$ws.Cells.Item(46, 5).Value = '  +1.37%  '
$ws.Cells.Item(47, 4).Value = '2.87'
$ws.Cells.Item(47, 5).Value = '  +3.23%  '
$ws.Cells.Item(48, 4).Value = '56.07'
$ws.Cells.Item(48, 5).Value = '  +8.05%  '
$ws.Cells.Item(49, 4).Value = '2.538.73'
$ws.Cells.Item(50, 5).Value = '  +5.64%  '
$ws.Cells.Item(51, 5).Value = '  +1.48%  '
